$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell C2 with value "ENGLISH"
$ws.Range("C2").Value = "ENGLISH"

# Clear the promotion validity cell R2 (previously "Promotion valid until  31th Dec 2021")
$ws.Range("R2").ClearContents()
